$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on all target cells first, so that numeric-looking
# strings (e.g. "0.999", "59.43") are stored as text, matching the source data,
# and not auto-converted to numbers by Excel.
$cellRefs = @(
"D2", "E2", "D3", "E3", "E4", "E5", "D6", "E6", "D7", "E7", "E8", "E9", "E10", "E11", "B12", "C12", "D12", "E12", "B13", "C13", "D13", "E13", "D14", "E14", "D15", "E15", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "D23", "E23", "E24", "E25", "E26", "D27", "D28", "E28", "D29", "E29", "E30", "E31", "D32", "E32", "E33", "D34", "E34", "D35", "E35", "D36", "E36", "E37", "D38", "E38", "E39", "D40", "E41", "D42", "E42", "D43", "E43", "E44", "D45", "E45", "E46", "E47", "E48", "E49", "E50", "D51", "E51"
)
foreach ($ref in $cellRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = "37.891.42"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "2.029.40"
$ws.Range("E3").Value = "  -1.08%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("E5").Value = "  -1.02%  "
$ws.Range("D6").Value = "0.609"
$ws.Range("E6").Value = "  -1.08%  "
$ws.Range("D7").Value = "59.43"
$ws.Range("E7").Value = "  +2.18%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  -0.73%  "
$ws.Range("E10").Value = "  +0.25%  "
$ws.Range("E11").Value = "  +0.56%  "
$ws.Range("B12").Value = "Chainlink"
$ws.Range("C12").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D12").Value = "14.56"
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "2.330.19"
$ws.Range("E13").Value = "  -1.10%  "
$ws.Range("D14").Value = "21.12"
$ws.Range("E14").Value = "  +2.06%  "
$ws.Range("D15").Value = "0.759"
$ws.Range("E15").Value = "  +1.22%  "
$ws.Range("E16").Value = "  -1.87%  "
$ws.Range("D17").Value = "2.029.15"
$ws.Range("E17").Value = "  -0.86%  "
$ws.Range("D18").Value = "37.805.69"
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("D19").Value = "6.02"
$ws.Range("E19").Value = "  -1.85%  "
$ws.Range("D20").Value = "70.00"
$ws.Range("E20").Value = "  +0.49%  "
$ws.Range("D21").Value = "0.0₃0823"
$ws.Range("E21").Value = "  -1.03%  "
$ws.Range("D22").Value = "224.90"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("E24").Value = "  -1.66%  "
$ws.Range("E25").Value = "  -2.03%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").Value = "165.25"
$ws.Range("D28").Value = "0.129"
$ws.Range("E28").Value = "  -2.53%  "
$ws.Range("D29").Value = "18.88"
$ws.Range("E29").Value = "  -0.77%  "
$ws.Range("E30").Value = "  -4.47%  "
$ws.Range("E31").Value = "  +0.65%  "
$ws.Range("D32").Value = "4.43"
$ws.Range("E32").Value = "  -2.29%  "
$ws.Range("E33").Value = "  +0.26%  "
$ws.Range("D34").Value = "4.51"
$ws.Range("E34").Value = "  -1.38%  "
$ws.Range("D35").Value = "0.0603"
$ws.Range("E35").Value = "  -1.46%  "
$ws.Range("D36").Value = "6.37"
$ws.Range("E36").Value = "  +6.73%  "
$ws.Range("E37").Value = "  -2.89%  "
$ws.Range("D38").Value = "3.24"
$ws.Range("E38").Value = "  -2.41%  "
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("D40").Value = "1.518.35"
$ws.Range("E41").Value = "  +0.22%  "
$ws.Range("D42").Value = "96.59"
$ws.Range("E42").Value = "  -1.11%  "
$ws.Range("D43").Value = "16.75"
$ws.Range("E43").Value = "  +0.91%  "
$ws.Range("E44").Value = "  -0.57%  "
$ws.Range("D45").Value = "0.0916"
$ws.Range("E45").Value = "  -1.84%  "
$ws.Range("E46").Value = "  -1.85%  "
$ws.Range("E47").Value = "  -3.56%  "
$ws.Range("E48").Value = "  -0.63%  "
$ws.Range("E49").Value = "  -0.29%  "
$ws.Range("E50").Value = "  +0.99%  "
$ws.Range("D51").Value = "2.217.97"
$ws.Range("E51").Value = "  -1.17%  "

# Restore default (Normal) style so no stray formatting/style id is introduced
foreach ($ref in $cellRefs) {
    $ws.Range($ref).Style = "Normal"
}
